$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new Price (column D) / new Volume(1h) (column E) text values.
# Only rows/columns actually present in the update are listed.
$updates = @(
    @{ Row = 2;  D = "42.839.36";  E = "  -0.08%  " },
    @{ Row = 3;  D = "2.569.26";   E = "  +1.41%  " },
    @{ Row = 4;                   E = "  +0.01%  " },
    @{ Row = 5;  D = "314.07";     E = "  -0.60%  " },
    @{ Row = 6;  D = "99.38";      E = "  +3.42%  " },
    @{ Row = 7;                   E = "  -0.30%  " },
    @{ Row = 9;  D = "0.533";      E = "  -0.09%  " },
    @{ Row = 10; D = "36.00";      E = "  -0.53%  " },
    @{ Row = 11; D = "0.0810";     E = "  -0.15%  " },
    @{ Row = 12; D = "7.48";       E = "  -1.20%  " },
    @{ Row = 13; D = "2.965.61";   E = "  +1.48%  " },
    @{ Row = 14;                  E = "  -0.84%  " },
    @{ Row = 15; D = "15.76";      E = "  +3.37%  " },
    @{ Row = 16; D = "2.499.04";   E = "  -2.17%  " },
    @{ Row = 17; D = "0.847";      E = "  -1.04%  " },
    @{ Row = 18; D = "42.898.48";  E = "  -0.07%  " },
    @{ Row = 19; D = "6.78";       E = "  -0.04%  " },
    @{ Row = 20; D = "12.57";      E = "  -2.41%  " },
    @{ Row = 21; D = "0.0₃0963";   E = "  -0.14%  " },
    @{ Row = 22; D = "69.42";      E = "  -0.77%  " },
    @{ Row = 23; D = "249.35";     E = "  -1.87%  " },
    @{ Row = 24; D = "2.95";       E = "  +0.11%  " },
    @{ Row = 25;                  E = "  -0.36%  " },
    @{ Row = 26; D = "26.96";      E = "  +0.63%  " },
    @{ Row = 27; D = "0.994";      E = "  -0.61%  " },
    @{ Row = 28; D = "2.36";       E = "  -2.66%  " },
    @{ Row = 29; D = "40.51";      E = "  -0.68%  " },
    @{ Row = 30; D = "10.26";      E = "  -1.14%  " },
    @{ Row = 31; D = "157.41";     E = "  -0.16%  " },
    @{ Row = 32; D = "5.80";       E = "  -2.10%  " },
    @{ Row = 33; D = "3.38";       E = "  +1.03%  " },
    @{ Row = 34;                  E = "  -3.16%  " },
    @{ Row = 35; D = "0.0801";     E = "  +2.39%  " },
    @{ Row = 36;                  E = "  +0.27%  " },
    @{ Row = 37; D = "18.81";      E = "  -2.09%  " },
    @{ Row = 38;                  E = "  +10.04%  " },
    @{ Row = 39;                  E = "  +0.05%  " },
    @{ Row = 40;                  E = "  -0.15%  " },
    @{ Row = 41; D = "23.39";      E = "  -0.43%  " },
    @{ Row = 42; D = "4.10";       E = "  +6.45%  " },
    @{ Row = 43;                  E = "  -0.11%  " },
    @{ Row = 44;                  E = "  -0.62%  " },
    @{ Row = 45; D = "3.24";       E = "  -2.50%  " },
    @{ Row = 46; D = "2.008.26";   E = "  -1.67%  " },
    @{ Row = 47;                  E = "  -1.24%  " },
    @{ Row = 48; D = "2.816.61";   E = "  +1.40%  " },
    @{ Row = 49;                  E = "  +2.25%  " },
    @{ Row = 50; D = "74.69";      E = "  -0.76%  " },
    @{ Row = 51; D = "81.61";      E = "  -4.02%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        # Prefix with an apostrophe so Excel stores the numeric-looking
        # text (e.g. "36.00", "0.0810") as literal text instead of
        # silently re-typing it as a number. Reapplying the "Normal"
        # style afterwards clears the quote-prefix formatting flag so
        # the cell keeps its original (default) style.
        $ws.Range("D$r").Value = "'" + $u.D
        $ws.Range("D$r").Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$r").Value = $u.E
    }
}
